$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 2) under the existing header row.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.6659568311637312
$ws.Range("C2").Value = 2.215936028923374
$ws.Range("D2").Value = 1.170718510326486

# A2 picks up the same header style (bold, centered, bordered) as B1:D1.
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
